# Refresh the cryptocurrency price/volume data (and two rank swaps) to match
# the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '54.134.94'
$ws.Range("E2").Value = '  -0.95%  '

# Row 3
$ws.Range("D3").Value = '2.270.27'
$ws.Range("E3").Value = '  -1.09%  '

# Row 4
$ws.Range("E4").Value = '  +0.60%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '498.12'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.59%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '128.45'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.93%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.42%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.526'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.51%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0950'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.06%  '

# Row 10
$ws.Range("E10").Value = '  +0.39%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.334'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +2.96%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.70'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +1.73%  '

# Row 13
$ws.Range("D13").Value = '2.672.73'
$ws.Range("E13").Value = '  -0.74%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '22.61'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +4.04%  '

# Row 15
$ws.Range("D15").Value = '54.129.48'

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000129'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.07%  '

# Row 17
$ws.Range("D17").Value = '2.280.14'
$ws.Range("E17").Value = '  -0.81%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.20'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +1.87%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.12'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +1.52%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '302.69'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.78%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.31'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -2.11%  '

# Row 22
$ws.Range("E22").Value = '  +0.08%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '61.04'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -2.98%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.00'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.09%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.149'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.94%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.28'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +2.36%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '170.61'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.64%  '

# Row 28
$ws.Range("E28").Value = '  +0.61%  '

# Row 29
$ws.Range("B29").Value = 'Aptos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.90'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.23%  '

# Row 30
$ws.Range("B30").Value = 'PEPE'
$ws.Range("C30").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D30").Value = '0.0₃0681'
$ws.Range("E30").Value = '  -0.55%  '

# Row 31
$ws.Range("E31").Value = '  +0.21%  '

# Row 32
$ws.Range("E32").Value = '  -0.03%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '17.71'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.62%  '

# Row 34
$ws.Range("B34").Value = 'FirstDigitalUSD'
$ws.Range("C34").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.998'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.04%  '

# Row 35
$ws.Range("B35").Value = 'SuiNetwork'
$ws.Range("C35").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.952'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +9.95%  '

# Row 36
$ws.Range("E36").Value = '  -1.22%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.68'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.79%  '

# Row 38
$ws.Range("E38").Value = '  -0.68%  '

# Row 39
$ws.Range("E39").Value = '  -0.24%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.35'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +0.40%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.81'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -1.51%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '124.77'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -3.06%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0490'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +1.48%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0890'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.55%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.544'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -1.52%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '238.27'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -2.05%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.371'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.72%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0204'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.62%  '

# Row 49
$ws.Range("E49").Value = '  +0.54%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '16.13'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -1.64%  '

# Row 51
$ws.Range("E51").Value = '  -0.50%  '
